{"js": "// 1. Remove the existing \"_GoBack\" bookmark (currently sitting in the empty\n//    centered paragraph right before the \"\u0421\u0441\u044b\u043b\u043a\u0430 \u043d\u0430 \u0440\u0435\u043f\u043e\u0437\u0438\u0442\u043e\u0440\u0438\u0439:\" heading).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Update the repository link text from the \"09. Introduction to IoT\"\n//    folder to the \"10. Flexible methodologies (...)\" folder, while keeping\n//    the existing run (and its formatting / rsid attributes) intact.\nconst oldUrl =\n  \"https://github.com/olgashenkel/GeekBrains-specialization-ELECTIVES/tree/main/09.%20Introduction%20to%20IoT\";\nconst newUrl =\n  \"https://github.com/olgashenkel/GeekBrains-specialization-ELECTIVES/tree/main/10.%20Flexible%20methodologies%20(Agile%2C%20SCRUM%2C%20Kanban%20and%20others)\";\n\nconst body = context.document.body;\nconst results = body.search(oldUrl, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the repository URL text to replace.\");\n}\n\nconst urlRange = results.items[0];\nurlRange.insertText(newUrl, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Re-insert the \"_GoBack\" bookmark, now collapsed at the very end of the\n//    paragraph that holds the (updated) link, i.e. right after the trailing\n//    space run and before the paragraph mark.\nconst paragraph = urlRange.paragraphs.getFirst();\nconst endRange = paragraph.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldUrl = \"https://github.com/olgashenkel/GeekBrains-specialization-ELECTIVES/tree/main/09.%20Introduction%20to%20IoT\"\n$newUrl = \"https://github.com/olgashenkel/GeekBrains-specialization-ELECTIVES/tree/main/10.%20Flexible%20methodologies%20(Agile%2C%20SCRUM%2C%20Kanban%20and%20others)\"\n$newUrlMarker = \"Flexible%20methodologies\"\n\n# 1. Remove the existing \"_GoBack\" bookmark (currently sitting in the empty\n#    centered paragraph right before the \"\u0421\u0441\u044b\u043b\u043a\u0430 \u043d\u0430 \u0440\u0435\u043f\u043e\u0437\u0438\u0442\u043e\u0440\u0438\u0439:\" heading).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Update the repository link text from the \"09. Introduction to IoT\"\n#    folder to the \"10. Flexible methodologies (...)\" folder.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute($oldUrl, $false, $false, $false, $false, $false, $true, 1, $false, $newUrl, 2) | Out-Null\n\n# 3. Re-insert the \"_GoBack\" bookmark, now collapsed at the very end of the\n#    paragraph that holds the (updated) link, i.e. right after the trailing\n#    space run and before the paragraph mark.\n#    A temporary marker character is used because collapsing a range exactly\n#    at end-of-paragraph-content otherwise anchors the bookmark to the start\n#    of the paragraph instead of the intended position.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Contains($newUrlMarker)) {\n        $insertPos = $p.Range.End - 1\n        $marker = $p.Range.Duplicate\n        $marker.Start = $insertPos\n        $marker.End = $insertPos\n        $marker.InsertAfter(\"@\")\n        break\n    }\n}\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Contains($newUrlMarker)) {\n        $bmPos = $p.Range.End - 2\n        $bmRange = $p.Range.Duplicate\n        $bmRange.Start = $bmPos\n        $bmRange.End = $bmPos\n        $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n        break\n    }\n}\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Contains($newUrlMarker)) {\n        $markerStart = $p.Range.End - 2\n        $markerRange = $p.Range.Duplicate\n        $markerRange.Start = $markerStart\n        $markerRange.End = $markerStart + 1\n        $markerRange.Delete()\n        break\n    }\n}\n"}
